$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.400.80"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.70%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.839.81"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "261.23"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -4.05%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5188"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.66%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3263"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06777"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.59"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -6.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7761"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07756"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.844.02"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "87.56"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.25%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007967"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.425.45"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.068.56"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.610"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.505"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.965"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.78"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.179"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -7.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.652"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.92"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.31"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.160"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.107"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08688"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04816"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7208"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.88%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.847"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.14%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01773"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.216"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4818"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -5.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9076"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "111.09"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.33%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.698"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05920"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4150"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -5.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.999"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.96"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -7.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8850"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.90%  "
